# add property for npc
# Adds a new "Height" column (AB) to the NPC sheet with a value of 2 for
# every data row (rows 2-21), mirroring column AA's formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("AB1").Value = "Height"

# Data rows: Height = 2 for every NPC row
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 28).Value = 2
}

# Match column AB's width to column AA's width
$ws.Columns.Item(28).ColumnWidth = $ws.Columns.Item(27).ColumnWidth

# Reflect the author's final selection on the new column
$ws.Range("AB2:AB21").Select() | Out-Null
